# Update proforma "Inputs and Outputs" sheet title to "Results Summary and Inputs"
# and make that sheet the active/selected sheet (it was previously
# "Developer Cash Flow" that was active).

$wb = $excel.ActiveWorkbook

$wsInputs = $wb.Worksheets.Item("Inputs and Outputs")
$wsInputs.Range("A1").Value = "Results Summary and Inputs"

# Reset selection on the Inputs and Outputs sheet to A2 (was D2 before)
$wsInputs.Range("A2").Select()

# Make "Inputs and Outputs" the active sheet/tab (previously "Developer Cash Flow" was active)
$wsInputs.Activate()
